$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new day's data (date serial 46001 = 2025-12-10) for both
# stations, reusing the existing "四方坪站" / "高岭站" shared strings.
$ws.Range("A20").Value = 46001
$ws.Range("B20").Value = "四方坪站"
$ws.Range("C20").Value = 8498.33
$ws.Range("D20").Value = 7363
$ws.Range("E20").Value = 2770.04
$ws.Range("F20").Value = 374

$ws.Range("A21").Value = 46001
$ws.Range("B21").Value = "高岭站"
$ws.Range("C21").Value = 4483.63
$ws.Range("D21").Value = 3762.52
$ws.Range("E21").Value = 1194.08
$ws.Range("F21").Value = 164

# Match the author's final selection (I18) as recorded in the sheet view.
$ws.Range("I18").Select()
